$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Value" column for several system config rows
$ws.Range("B2").Value = "chrome"
$ws.Range("B3").Value = "172.41.46.23"
$ws.Range("B4").Value = "cs01"
$ws.Range("B9").Value = "Playwright_POC"
$ws.Range("B10").Value = "`nPP_2_6_3_B2_P2"

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("F11").Select()
